$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new data row at row 11 (new item #5), pushing the existing
# "total" row (old row 11) and the footer row (old row 12) down to rows
# 12 and 13 respectively.
# ---------------------------------------------------------------------------
$ws.Rows("11:11").Insert()

# Copy the formatting of the previous item row (row 10, the last product
# row) down into the freshly inserted row 11 so every cell picks up the
# same styles used by the other item rows (7-10).
$ws.Range("A10:Q10").Copy()
$ws.Range("A11:Q11").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Rows("11:11").RowHeight = 25.5

# Recreate the merged cell groups for the new row (matching the layout used
# by rows 7-10: A:B, C:G, H:K, L:M, N:O, with P/Q left as separate cells).
$ws.Range("A11:B11").Merge()
$ws.Range("C11:G11").Merge()
$ws.Range("H11:K11").Merge()
$ws.Range("L11:M11").Merge()
$ws.Range("N11:O11").Merge()

# ---------------------------------------------------------------------------
# Populate the new item row (#5) with its data.
# ---------------------------------------------------------------------------
# Item number
$ws.Range("A11").Value = 5

# Item name (plain text column already formatted as @)
$ws.Range("C11").Value = "راجون احمر 250 مل"

# Current balance (plain text column already formatted as @)
$ws.Range("H11").Value = "7:0"

# Order limit - column uses a numeric display format, so switch it to text
# just long enough to store the literal string, then restore the original
# numeric format (this keeps the cell's style id identical to the sibling
# rows while still writing a text value, matching the source workbook).
$fmtL = $ws.Range("L11").NumberFormat
$ws.Range("L11").NumberFormat = "@"
$ws.Range("L11").Value = "0"
$ws.Range("L11").NumberFormat = $fmtL

# Price (plain text column already formatted as @)
$ws.Range("N11").Value = "35.00"

# Sale price - numeric display format column, use the same text trick.
$fmtP = $ws.Range("P11").NumberFormat
$ws.Range("P11").NumberFormat = "@"
$ws.Range("P11").Value = "35.0000"
$ws.Range("P11").NumberFormat = $fmtP

# Number of transactions (plain text column already formatted as @)
$ws.Range("Q11").Value = "1:0"

# ---------------------------------------------------------------------------
# Update the total row (now row 12): add the new item's sale price to the
# previous total (108 + 35 = 143).
# ---------------------------------------------------------------------------
$ws.Range("P12").Value = 143

# ---------------------------------------------------------------------------
# Update the footer timestamp (now row 13) to reflect the new export time.
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = "Friday, 3 October, 2025 11:41 AM"

Write-Host "Applied update: inserted item #5 (راجون احمر 250 مل) and refreshed totals/timestamp."
